$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet has 14 data rows (rows 2-15): line1..line6 (rows 2-7) and
# extr1..extr8 (rows 8-15). Two new entries ("line7", "line8") are inserted
# into the dataset right after line6 and before extr1, so extr1..extr8 shift
# down by two rows (to rows 10-17), and two brand-new rows (16-17) are used
# for the last two extr entries that no longer fit in the old row range.

# Step 1: shift existing rows 8-15 (extr1..extr8) down to rows 10-17.
# Copy from the bottom up so source rows aren't overwritten before being read.
for ($r = 15; $r -ge 8; $r--) {
    $destRow = $r + 2
    $ws.Range("A" + $r + ":E" + $r).Copy($ws.Range("A" + $destRow))
}

# Step 2: write the new "line7" / "line8" rows into the now-vacant rows 8-9.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $false

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Step 3: update the shifted extr1..extr8 rows (10-17) with their final values.
# Column A is the sequential record index (0-based), columns C/D are the new
# endpoint values, and column E is the in_service flag.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true
